$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("techdata")

# The "level" column (D) used the value "biddingzone" for every generator
# row; rename it to "node" everywhere it appears in that column.
$ws.Columns("D").Replace("biddingzone", "node", 1)

# Restore the author's on-screen selection/viewport for the sheet.
$ws.Range("E72").Select()
